# Dic_Einheit.xlsx update:
#   - insert a new entry "E_GINI" (Gini-Koeffizient / Gini coefficient)
#     right after E_BNEUR (new row 5), shifting the rest of the table down
#   - rename the bottom two entries IDX_1990 / IDX_2030 to
#     E_IDX_1990 / E_IDX_2030 and move them up so the list stays sorted
#     (right after E_HAPD)
#   - E_TEUR ends up as the new last row (22)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank row at row 5 so the table grows from 21 to 22 data
#    rows (dimension A1:C21 -> A1:C22) and every following row shifts
#    down by one.
$ws.Rows("5:5").Insert()

# Give the freshly inserted row the same formatting as the rest of the
# data rows (style "s=4" in the original file) by copying it from a
# neighbouring data row.
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)

# 2) Write the final, reordered content for every data row (2-22).
$ws.Range("A2").Value = "E_10H3"
$ws.Range("B2").Value = "1 000"
# "1.000" looks like a number to Excel's auto-detection, so force the
# cell to Text before writing it, then restore the normal data-row
# formatting (style) from a neighbouring cell.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1.000"
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

$ws.Range("A3").Value = "E_10H7"
$ws.Range("B3").Value = "Millionen"
$ws.Range("C3").Value = "Millions"

$ws.Range("A4").Value = "E_BNEUR"
$ws.Range("B4").Value = "Milliarden EUR"
$ws.Range("C4").Value = "Billion EUR"

$ws.Range("A5").Value = "E_GINI"
$ws.Range("B5").Value = "Gini-Koeffizient"
$ws.Range("C5").Value = "Gini coefficient"

$ws.Range("A6").Value = "E_HAPD"
$ws.Range("B6").Value = "Hektar pro Tag"
$ws.Range("C6").Value = "Hectre per day"

$ws.Range("A7").Value = "E_IDX_1990"
$ws.Range("B7").Value = "1990 = 100"
$ws.Range("C7").Value = "1990 = 100"

$ws.Range("A8").Value = "E_IDX_2030"
$ws.Range("B8").Value = "2030 = 100"
$ws.Range("C8").Value = "2030 = 100"

$ws.Range("A9").Value = "E_IDX2000"
$ws.Range("B9").Value = "2000 = 100"
$ws.Range("C9").Value = "2000 = 100"

$ws.Range("A10").Value = "E_IDX2005"
$ws.Range("B10").Value = "2005 = 100"
$ws.Range("C10").Value = "2005 = 100"

$ws.Range("A11").Value = "E_IDX2008"
$ws.Range("B11").Value = "2008 = 100"
$ws.Range("C11").Value = "2008 = 100"

$ws.Range("A12").Value = "E_KGPHA"
$ws.Range("B12").Value = "Kilogramm pro Hektar"
$ws.Range("C12").Value = "Kilogram per hectare"

$ws.Range("A13").Value = "E_MILIGPL"
$ws.Range("B13").Value = "Miligramm pro Liter"
$ws.Range("C13").Value = "Miligrams per litre"

$ws.Range("A14").Value = "E_MIN"
$ws.Range("B14").Value = "Minuten"
$ws.Range("C14").Value = "Minutes"

$ws.Range("A15").Value = "E_MNEUR"
$ws.Range("B15").Value = "Millionen EUR"
$ws.Range("C15").Value = "Million EUR"

$ws.Range("A16").Value = "E_MNEW"
$ws.Range("B16").Value = "Millionen Einwohner/ -innen"
$ws.Range("C16").Value = "Million inhabitants"

$ws.Range("A17").Value = "E_NUM"
$ws.Range("B17").Value = "Anzahl"
$ws.Range("C17").Value = "Number"

$ws.Range("A18").Value = "E_P10H6EWN"
$ws.Range("B18").Value = "Je 100 000 Einwohner/ -innen"
$ws.Range("C18").Value = "Per 100,000 inhabitants"

$ws.Range("A19").Value = "E_P10H6EWNU70"
$ws.Range("B19").Value = "Je 100 000 Einwohner/ -innen unter 70 Jahren (ohne unter 1-Jährige)"
$ws.Range("C19").Value = "Per 100,000 inhabitants under 70 years (excluding under 1 year olds)"

$ws.Range("A20").Value = "E_PRZNT"
$ws.Range("B20").Value = "Prozent"
$ws.Range("C20").Value = "Percentage"

$ws.Range("A21").Value = "E_QMPA"
$ws.Range("B21").Value = "m² pro Jahr"
$ws.Range("C21").Value = "m² per year"

$ws.Range("A22").Value = "E_TEUR"
$ws.Range("B22").Value = "1 000 EUR"
$ws.Range("C22").Value = "1.000 EUR"
